# Append a new data row (row 95) to each of the 4 worksheets, mirroring the
# structure of the existing rows (time stamp, hex-string fields, numeric
# decoded fields).

$wb = $excel.ActiveWorkbook

# Data for the new row on each worksheet, in column order A..I
# (scientific-notation literals aren't accepted by the parser, so the
# G-column values are parsed from strings into doubles instead)
$gLft1 = [double]"7.598631275147109e+23"
$gLft2 = [double]"5.68432987514711e+23"
$gPlt1 = [double]"5.68631262647114e+23"
$gPlt2 = [double]"9.85046333984776e+23"

$rowsData = @{
    "FE_LFT_#1" = @(45881.49295138889, "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x01,0x04", "0xf", 380, $gLft1, 260, 15)
    "FE_LFT_#2" = @(45881.49295138889, "0x01,0x90", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x18", "0xe", 400, $gLft2, 280, 14)
    "FE_PLT_#1" = @(45881.49295138889, "0x00,0x6e", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x5C", "0x3", 110, $gPlt1, 92, 3)
    "FE_PLT_#2" = @(45881.49295138889, "0x00,0x6e", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x5B", "0x3", 110, $gPlt2, 91, 3)
}

foreach ($sheetName in $rowsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $rowsData[$sheetName]

    $newRow = 95

    # Column A: numeric date/time value, formatted like the rows above it.
    $cellA = $ws.Cells.Item($newRow, 1)
    $cellA.Value = $values[0]
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Columns B-E: hex strings stored as text.
    $ws.Cells.Item($newRow, 2).Value = $values[1]
    $ws.Cells.Item($newRow, 3).Value = $values[2]
    $ws.Cells.Item($newRow, 4).Value = $values[3]
    $ws.Cells.Item($newRow, 5).Value = $values[4]

    # Columns F-I: numeric values.
    $ws.Cells.Item($newRow, 6).Value = $values[5]
    $ws.Cells.Item($newRow, 7).Value = $values[6]
    $ws.Cells.Item($newRow, 8).Value = $values[7]
    $ws.Cells.Item($newRow, 9).Value = $values[8]
}
